$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 821.55554
$ws.Range("J17").Value = 621.375
$ws.Range("L17").Value = 1864.125
$ws.Range("N17").Value = -2200.125
$ws.Range("H74").Value = 6952.1816
$ws.Range("I74").Value = 7147.45
$ws.Range("K74").Value = 7147.45
$ws.Range("M74").Value = -6211.45
$ws.Range("H77").Value = 6952.1816
$ws.Range("I77").Value = 7147.45
$ws.Range("K77").Value = 35737.25
$ws.Range("M77").Value = -31057.25
$ws.Range("H98").Value = 2332.476
$ws.Range("I98").Value = 2432.7778
$ws.Range("J98").Value = 1730.6666
$ws.Range("K98").Value = 2432.7778
$ws.Range("L98").Value = 1730.6666
$ws.Range("M98").Value = -934.7777999999998
$ws.Range("N98").Value = -4726.6666
$ws.Range("H121").Value = 4849.5
$ws.Range("J121").Value = 4849.5
$ws.Range("L121").Value = 14548.5
$ws.Range("N121").Value = -18042.5
$ws.Range("H122").Value = 2332.476
$ws.Range("I122").Value = 2432.7778
$ws.Range("J122").Value = 1730.6666
$ws.Range("K122").Value = 7298.3334
$ws.Range("L122").Value = 5191.9998
$ws.Range("M122").Value = -4848.3334
$ws.Range("N122").Value = -10091.9998
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H137").Value = 17246470
$ws.Range("I137").Value = 23811158
$ws.Range("J137").Value = 14166.5
$ws.Range("K137").Value = 71433474
$ws.Range("L137").Value = 42499.5
$ws.Range("M137").Value = -71430924
$ws.Range("N137").Value = -47599.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 824886.0600000001
$ws.Range("I74").Value = 947631.9
$ws.Range("K74").Value = 947631.9
$ws.Range("M74").Value = -946757.9
$ws.Range("H77").Value = 824886.0600000001
$ws.Range("I77").Value = 947631.9
$ws.Range("K77").Value = 4738159.5
$ws.Range("M77").Value = -4733791.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1477.5883
$ws.Range("I107").Value = 1170.2106
$ws.Range("K107").Value = 1170.2106
$ws.Range("M107").Value = 749.7893999999999
$ws.Range("H134").Value = 3473748.8
$ws.Range("I134").Value = 1601.8572
$ws.Range("J134").Value = 27778778
$ws.Range("K134").Value = 4805.571599999999
$ws.Range("L134").Value = 83336334
$ws.Range("M134").Value = -2270.571599999999
$ws.Range("N134").Value = -83341404

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8236484
$ws.Range("I31").Value = 2501264.5
$ws.Range("J31").Value = 100000000
$ws.Range("K31").Value = 2501264.5
$ws.Range("L31").Value = 100000000
$ws.Range("M31").Value = -2500969.5
$ws.Range("N31").Value = -100000590
$ws.Range("H34").Value = 8236484
$ws.Range("I34").Value = 2501264.5
$ws.Range("J34").Value = 100000000
$ws.Range("K34").Value = 2501264.5
$ws.Range("L34").Value = 100000000
$ws.Range("M34").Value = -2501062.5
$ws.Range("N34").Value = -100000404
$ws.Range("H58").Value = 6634822.5
$ws.Range("I58").Value = 13890629
$ws.Range("K58").Value = 13890629
$ws.Range("M58").Value = -13890426
$ws.Range("H99").Value = 2576.1428
$ws.Range("I99").Value = 2077.3333
$ws.Range("J99").Value = 2950.25
$ws.Range("K99").Value = 2077.3333
$ws.Range("L99").Value = 2950.25
$ws.Range("M99").Value = -579.3332999999998
$ws.Range("N99").Value = -5946.25
$ws.Range("H107").Value = 659.5238000000001
$ws.Range("I107").Value = 453.25
$ws.Range("K107").Value = 453.25
$ws.Range("M107").Value = 1466.75
$ws.Range("H126").Value = 2576.1428
$ws.Range("I126").Value = 2077.3333
$ws.Range("J126").Value = 2950.25
$ws.Range("K126").Value = 6231.999899999999
$ws.Range("L126").Value = 8850.75
$ws.Range("M126").Value = -3761.999899999999
$ws.Range("N126").Value = -13790.75
$ws.Range("H132").Value = 2616.5454
$ws.Range("I132").Value = 2489.1765
$ws.Range("K132").Value = 7467.529500000001
$ws.Range("M132").Value = -4937.529500000001
$ws.Range("H136").Value = 6634822.5
$ws.Range("I136").Value = 13890629
$ws.Range("K136").Value = 41671887
$ws.Range("M136").Value = -41669337

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 117.57895
$ws.Range("I2").Value = 24.636364
$ws.Range("K2").Value = 147.818184
$ws.Range("M2").Value = -34.818184
$ws.Range("H5").Value = 1843427
$ws.Range("I5").Value = 1553305.9
$ws.Range("K5").Value = 4659917.699999999
$ws.Range("M5").Value = -4659805.699999999
$ws.Range("H107").Value = 2965.5386
$ws.Range("J107").Value = 4564.5713
$ws.Range("L107").Value = 13693.7139
$ws.Range("N107").Value = -17533.7139
$ws.Range("H114").Value = 1762.6666
$ws.Range("I114").Value = 216.5
$ws.Range("J114").Value = 2999.6
$ws.Range("K114").Value = 649.5
$ws.Range("L114").Value = 8998.799999999999
$ws.Range("M114").Value = 2604.5
$ws.Range("N114").Value = -15506.8
$ws.Range("H131").Value = 3299.0425
$ws.Range("I131").Value = 548.625
$ws.Range("J131").Value = 4718.613
$ws.Range("K131").Value = 1645.875
$ws.Range("L131").Value = 14155.839
$ws.Range("M131").Value = 3394.125
$ws.Range("N131").Value = -24235.839
$ws.Range("H134").Value = 5338.125
$ws.Range("I134").Value = 2228.611
$ws.Range("K134").Value = 6685.833
$ws.Range("M134").Value = -1615.833
$ws.Range("H135").Value = 1843427
$ws.Range("I135").Value = 1553305.9
$ws.Range("K135").Value = 13979753.1
$ws.Range("M135").Value = -13977218.1
$ws.Range("H139").Value = 6755.5386
$ws.Range("I139").Value = 3924.6667
$ws.Range("K139").Value = 11774.0001
$ws.Range("M139").Value = -6634.000100000001
$ws.Range("H140").Value = 3956.6
$ws.Range("I140").Value = 3396.2222
$ws.Range("K140").Value = 10188.6666
$ws.Range("M140").Value = -5008.6666

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18794
$ws.Range("I70").Value = 41177.4
$ws.Range("K70").Value = 41177.4
$ws.Range("M70").Value = -40907.4
$ws.Range("H73").Value = 18794
$ws.Range("I73").Value = 41177.4
$ws.Range("K73").Value = 41177.4
$ws.Range("M73").Value = -40241.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4046.8572
$ws.Range("I46").Value = 1249.5
$ws.Range("J46").Value = 5165.8
$ws.Range("K46").Value = 1249.5
$ws.Range("L46").Value = 5165.8
$ws.Range("M46").Value = -1061.5
$ws.Range("N46").Value = -5541.8
$ws.Range("H68").Value = 3660.3572
$ws.Range("I68").Value = 2195
$ws.Range("J68").Value = 10401
$ws.Range("K68").Value = 2195
$ws.Range("L68").Value = 10401
$ws.Range("M68").Value = -1446
$ws.Range("N68").Value = -11899
$ws.Range("H71").Value = 3660.3572
$ws.Range("I71").Value = 2195
$ws.Range("J71").Value = 10401
$ws.Range("K71").Value = 10975
$ws.Range("L71").Value = 52005
$ws.Range("M71").Value = -7231
$ws.Range("N71").Value = -59493
$ws.Range("H132").Value = 4389382.5
$ws.Range("I132").Value = 9262481
$ws.Range("J132").Value = 3594.2
$ws.Range("K132").Value = 27787443
$ws.Range("L132").Value = 10782.6
$ws.Range("M132").Value = -27784913
$ws.Range("N132").Value = -15842.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 69599.5
$ws.Range("J46").Value = 69599.5
$ws.Range("L46").Value = 69599.5
$ws.Range("N46").Value = -70061.5
$ws.Range("H122").Value = 42992.85
$ws.Range("I122").Value = 1249.625
$ws.Range("K122").Value = 3748.875
$ws.Range("M122").Value = -1298.875
$ws.Range("H126").Value = 1678.8667
$ws.Range("I126").Value = 1224.0834
$ws.Range("K126").Value = 3672.2502
$ws.Range("M126").Value = -1202.2502
$ws.Range("H132").Value = 83335580
$ws.Range("I132").Value = 166666670
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 500000010
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -499997480
$ws.Range("N132").Value = -18560
$ws.Range("H134").Value = 69599.5
$ws.Range("J134").Value = 69599.5
$ws.Range("L134").Value = 208798.5
$ws.Range("N134").Value = -213868.5
$ws.Range("H136").Value = 13248469
$ws.Range("I136").Value = 6212535.5
$ws.Range("J136").Value = 62500000
$ws.Range("K136").Value = 18637606.5
$ws.Range("L136").Value = 187500000
$ws.Range("M136").Value = -18635056.5
$ws.Range("N136").Value = -187505100
